$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for "Introduction to Deep Learning in Python" above the
# existing "Introduction to Network Analysis in Python" row (old row 30,
# new row 31), matching the style/format of the row it's being inserted in
# front of (inherits formatting from the row above automatically).
$ws.Rows.Item(30).Insert()
$ws.Cells.Item(30, 1).Value = "Introduction to Deep Learning in Python"
$ws.Cells.Item(30, 3).Value = 4

# Insert a new row for "Advanced Deep Learning with Keras" above the
# existing "Supervised Learning with scikit-learn" row (now at row 35
# after the previous insertion).
$ws.Rows.Item(35).Insert()
$ws.Cells.Item(35, 1).Value = "Advanced Deep Learning with Keras"

# Restore the selection to match where the author ended up after editing.
$ws.Range("C31").Select()
